# Generate Report for Handback
# The f1934dc0-44fe-469d-81ce-6c5a14bf8cf2 file has now been successfully
# handed back (in sync with en-US) for both the zh-cn and de-de locales, so
# update the Overview sheet and the per-locale report sheets accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the f1934dc0... file, zh-cn (E) and de-de (F) columns ---
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: row 3 is the f1934dc0... file ---
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "2016-08-26 16:49:15"
$zhcn.Range("P3").Value = "'"
$zhcn.Range("P3").Style = "Normal"

# --- de-de sheet: row 3 is the f1934dc0... file ---
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "2016-08-26 16:49:22"
$dede.Range("P3").Value = "'"
$dede.Range("P3").Style = "Normal"
